# Update Excel upload method (level2 screen) - manucost.xlsx datatemplate
# Row 1 holds SQL column types, row 2 holds column names/labels.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: column data types ---
# A1 "int(11)" stays unchanged.
$ws.Range("B1").Value = "varchar(50)"
$ws.Range("C1").Value = "varchar(50)"
$ws.Range("D1").Value = "varchar(8)"
$ws.Range("E1").Value = "varchar(50)"
$ws.Range("F1").Value = "int(20)"
# G1 "varchar(50)" stays unchanged.
$ws.Range("H1").Value = "varchar(50)"

# --- Row 2: column labels (now localized Korean labels) ---
$ws.Range("A2").Value = "(id)입력x"
$ws.Range("B2").Value = "사업장ID"
$ws.Range("C2").Value = "코스트센터ID"
$ws.Range("D2").Value = "년월"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "금액"
$ws.Range("G2").Value = "버젼ID"
$ws.Range("H2").Value = "관리항목"
